$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Create new row 5 (2021年) by copying the formatting of row 4's label cell,
# then overwrite with the new label and data values.
$ws.Cells.Item(4, 1).Copy($ws.Cells.Item(5, 1))
$ws.Cells.Item(5, 1).Value = "2021年"

$ws.Cells.Item(5, 2).Value = 35.849
$ws.Cells.Item(5, 3).Value = 38.454
$ws.Cells.Item(5, 4).Value = 26.12
$ws.Cells.Item(5, 5).Value = 31.311
$ws.Cells.Item(5, 6).Value = 38.71
$ws.Cells.Item(5, 7).Value = 39.878
$ws.Cells.Item(5, 8).Value = 40.152
$ws.Cells.Item(5, 9).Value = 37.746
$ws.Cells.Item(5, 10).Value = 28.889
$ws.Cells.Item(5, 11).Value = 34.354
$ws.Cells.Item(5, 12).Value = 63.564
$ws.Cells.Item(5, 13).Value = 18.701
